# Update cryptocurrency price/volume data on the active sheet (Sheet1).
# Source data cells are stored as text (inline strings), so we force a
# text number format before writing so Excel doesn't reinterpret the
# values as numbers/percentages and strip significant trailing zeros.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = [ordered]@{
    "D2"  = "289.14"
    "E2"  = "0.88%"
    "E3"  = "1.67%"
    "D4"  = "5.255"
    "E4"  = "3.68%"
    "D5"  = "0.07059"
    "E5"  = "6.00%"
    "D6"  = "7.462"
    "E6"  = "1.60%"
    "D7"  = "3.569"
    "E7"  = "5.16%"
    "D8"  = "1.395"
    "E8"  = "1.61%"
    "D9"  = "0.9058"
    "E9"  = "-3.86%"
    "E10" = "3.37%"
    "D11" = "0.07621"
    "E11" = "15.20%"
    "D12" = "0.07793"
    "E12" = "2.95%"
    "D13" = "0.02911"
    "E13" = "-1.59%"
    "D14" = "0.09023"
    "E14" = "0.17%"
    "D15" = "0.001590"
    "E15" = "1.01%"
    "D16" = "0.0006516"
    "E16" = "0.73%"
    "D17" = "0.006176"
    "E17" = "-2.38%"
    "D18" = "3.486"
    "E18" = "1.04%"
    "E19" = "-0.89%"
    "D20" = "0.3233"
    "E20" = "0.56%"
    "E21" = "2.78%"
    "D22" = "4.004"
    "E22" = "-2.21%"
    "E23" = "2.83%"
    "D24" = "0.04517"
    "E24" = "0.55%"
    "D25" = "0.001210"
    "E25" = "2.25%"
    "D26" = "0.004168"
    "E26" = "-7.22%"
    "D27" = "0.0001168"
    "E27" = "-6.71%"
    "D28" = "0.0001667"
    "E28" = "2.80%"
    "D40" = "0.04387"
    "E40" = "4.41%"
    "D41" = "0.007012"
    "E41" = "4.01%"
    "D42" = "0.1252"
    "E42" = "-0.49%"
    "D43" = "0.002066"
    "E43" = "2.14%"
    "D44" = "0.01192"
    "E44" = "-3.45%"
    "D45" = "0.00005836"
    "E45" = "2.83%"
    "D47" = "0.01298"
    "E47" = "-0.88%"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
}
